$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1: bold font, thin border all around, centered horizontally, top vertically
$r1 = $ws.Range("B1")
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1        # xlContinuous
$r1.Borders.Weight = 2           # xlThin

# Copy B1's formatting onto A2 so both share the same cell style (s="1")
# instead of creating a second, near-duplicate style entry.
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
